$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update mfd_sampletype (K) and mfd_areatype (L) casing, and fill in
# mfd_hab1 (N) / mfd_hab2 (O) grid-centroid habitat values for each
# agricultural sample row (rows 2-11).
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 11).Value = "Soil"
    $ws.Cells.Item($r, 12).Value = "Natural"
    $ws.Cells.Item($r, 14).Value = "Grassland formations"
    $ws.Cells.Item($r, 15).Value = "Semi-natural tall-herb humid meadows"
}
